$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value of 45203 (2023-10-04) for
# every data row (rows 2 through 410). Update it to 45204 (2023-10-05).
for ($row = 2; $row -le 410; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
